# "Generate Report for handback"
#
# The handback-status workbook tracks, per locale sheet, the handoff/handback
# round-trip for each localized file. A new handback round just completed for
# the first file in each locale (row 2), so its "Correspond Handoff Datetime"
# (column D) and "Correspond Handback DateTime" (column G) need to be
# refreshed with the timestamps of this new round.
#
# Use Value2 (not Value) so the timestamp strings are stored verbatim as text
# (matching the existing cells), instead of being auto-coerced into a date
# serial number with a new number-format style.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value2 = "2016-01-11 06:01:11"
$wsZhCn.Range("G2").Value2 = "2016-01-11 06:02:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value2 = "2016-01-11 06:01:30"
$wsDeDe.Range("G2").Value2 = "2016-01-11 06:02:54"
